$p = $ppt.ActivePresentation

# --- Update the datetimeFigureOut field text across all existing slides ---
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "9/20/2022") {
                $tr.Text = "9/25/2022"
            }
        }
    }
}

# --- Add the new "References" slide (slide 11) using the Title and Content layout ---
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "References"

$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange
$tr.Text = "How do you play Planning Poker? planningpoker.com. 2020.`rFlowchart Maker and Online Diagram Software. Draw.io. 2005-2022. JGraph Ltd.`rSoftware Requirements Specification. Dr. Kirstie Hawkey. 2011. Dalhousie University.`r`rReading for next lecture: Pressman Ch 9-11, Appendix 1"

$tr.Font.NameComplexScript = "Calibri"

$tr.Paragraphs(1).ActionSettings.Item(1).Hyperlink.Address = "https://planningpoker.com/"
$tr.Paragraphs(2).ActionSettings.Item(1).Hyperlink.Address = "https://www.drawio.com/"
$tr.Paragraphs(3).ActionSettings.Item(1).Hyperlink.Address = "https://www.dal.ca/"

$tr.Paragraphs(5).Font.Italic = -1

$tf.AutoSize = 2
$tf.MarginLeft = 7.2
$tf.MarginRight = 7.2
$tf.MarginTop = 3.6
$tf.MarginBottom = 3.6
$tf.VerticalAnchor = 1
$tf.Orientation = 1
